# Add the new data rows (41-45) to the "Online" ORA Errors sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of (Date serial, Error Count) to append below the existing data.
$newRows = @(
    @(46021, 48),
    @(46027, 90),
    @(46020, 47),
    @(46017, 45),
    @(46024, 51)
)

$startRow = 41

# Copy the date formatting (number format / style) from the last existing
# data row (A40) down onto the new date cells (A41:A45) before filling in
# the values, so the new dates keep the same date display format.
$ws.Range("A40").Copy() | Out-Null
$lastRow = $startRow + $newRows.Length - 1
$ws.Range("A" + $startRow + ":A" + $lastRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Update the view so the newly added rows are shown/selected, matching
# how Excel leaves the sheet scrolled/selected after data entry.
$ws.Rows("44:44").Select() | Out-Null
$win = $wb.Windows.Item(1)
$win.ScrollRow = 37
$win.ScrollColumn = 1
